# Apply the edit described by the commit:
#  "Task Interface und BackendService erstellt, in List-K. Backend eingebunden,
#   auslesen des Taskarrays in Tabelle"
#
# Adds two new worksheets after "Tabelle1":
#   - "30.09."        with a note about the frontend interface/service work
#   - "Notizen Präsi"  with a small list of notes, and becomes the active sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet "30.09." right after Tabelle1
$sheet30 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheet30.Name = "30.09."

# New sheet "Notizen Präsi" right after "30.09."
$sheetNotizen = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheetNotizen.Name = "Notizen Präsi"

# Fill in the notes (order chosen to match shared-string insertion order)
$sheetNotizen.Range("A1").Value = "Backend starten"
$sheetNotizen.Range("A2").Value = "Frontend starten"

$sheet30.Range("A2").Value = "interface und service im frontend angelegt+"

$sheetNotizen.Range("A3").Value = "klären === und == ?"

# Restore/update selections on each sheet
$ws1.Range("E3").Select() | Out-Null
$sheet30.Range("A3").Select() | Out-Null
$sheetNotizen.Range("A4").Select() | Out-Null
